# Commit message: "just backup - importing between files"
# This replicates a sheet copy: "0010" is duplicated to a new sheet
# named "00101" appended at the end of the workbook, with two minor
# text corrections in the copy (an extra blank line removed in F2 and C16).

$wb = $excel.ActiveWorkbook

# Copy the existing "0010" sheet to the very end of the workbook.
$src = $wb.Worksheets.Item("0010")
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# The freshly-copied sheet is now the last sheet; rename it to "00101".
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "00101"

# Two cells in the copy differ from the source by a stray blank line;
# rewrite them with the corrected text.
$newSheet.Range("F2").Value = " Số dư`nBalance`n"
$newSheet.Range("C16").Value = " MA_GD:107082158|K194111578,`nHoc`nphiHK01,HK01,HK01,HK0O1,HKO1,I`nHK01/2020-2021`n50000][201 DL06@1@1[BPMENT-F`n"
